# Append " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document."), as three separate runs:
#   " (", "Changed main", ")"
#
# A plain Range.InsertAfter() call would get silently coalesced back
# into the neighbouring run (same empty run formatting => the engine,
# like real Word, normalizes adjacent identical-format runs on save).
# Wrapping the inserts in tracked-change markers and then accepting
# each revision individually (instead of Revisions.AcceptAll, which
# forces a full-document reflow/cleanup pass) keeps the three new runs
# distinct in the saved OOXML while leaving the rest of the document
# untouched.

$d = $word.ActiveDocument

$p = $d.Paragraphs(1)
$pRange = $p.Range

# Exclude the trailing paragraph mark so the insert lands at the end
# of the visible text, not at the start of the next paragraph.
$insertPoint = $d.Range($pRange.Start, $pRange.End - 1)

$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true

$insertPoint.Collapse(0)
$insertPoint.InsertAfter(" (")

$insertPoint.Collapse(0)
$insertPoint.InsertAfter("Changed main")

$insertPoint.Collapse(0)
$insertPoint.InsertAfter(")")

$d.TrackRevisions = $wasTracking

# Accept the three insertions one at a time (not AcceptAll) so only
# the edited paragraph is touched.
while ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item(1).Accept()
}
